$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# J2: "001" -> "002" (force text, keep leading zero, then drop the quote-prefix style)
$ws.Range("J2").Value = "'002"
$ws.Range("J2").Style = "Normal"

# N2: REPORT_DATE text update
$ws.Range("N2").Value = "2020-06-30 00:00:00"

# Numeric value updates
$ws.Range("O2").Value = 35393856.22
$ws.Range("P2").Value = 167929441.98
$ws.Range("Q2").Value = 129678107.83

# R2: cleared to an empty text cell (matches the blank inlineStr cells elsewhere in the row)
$ws.Range("R2").Value = "'"
$ws.Range("R2").Style = "Normal"

$ws.Range("S2").Value = 106817254.74
$ws.Range("T2").Value = 106817254.74

# U2: cleared to an empty text cell
$ws.Range("U2").Value = "'"
$ws.Range("U2").Style = "Normal"

$ws.Range("V2").Value = 2630741.43
$ws.Range("W2").Value = 7988377.15
$ws.Range("X2").Value = 391372.77
$ws.Range("Y2").Value = 37641750.21
$ws.Range("Z2").Value = 40284970.65
$ws.Range("AA2").Value = 4891114.43
$ws.Range("AG2").Value = 1408533.84
$ws.Range("AP2").Value = 49.348143374
$ws.Range("AQ2").Value = 13.128995951709
$ws.Range("AR2").Value = 22.55235770725
$ws.Range("AS2").Value = 32592943
$ws.Range("AT2").Value = 12.525505007785
